# Actualización automática 2025-11-26 10:30:08
# Updates sales figures for CASTRO ALCIVAR EDA MARIA / MAD&DECO S.A.
# across the three report sheets, propagating the new totals and
# recalculated "cumplimiento" (% achieved) figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales by product group, row 34 = client)
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D34").Value = 5280.57    # 240X80 PORCELANATO
$wsGrupo.Range("L34").Value = 506.88     # PIEDRA SINTERIZADA
$wsGrupo.Range("M34").Value = 2473.79    # PORCELANATO

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (monthly sales, row 34 = client, row 61 = total)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F34").Value = 7840.56   # noviembre, client row
$wsMensual.Range("F61").Value = 57481.96  # noviembre, total row

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (monthly compliance / achievement)
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 21657.75
$wsCumpl.Range("E3").Value = -6832.34
$wsCumpl.Range("F3").Value = 1.460853359198835

# Row 11 - PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 8165.68
$wsCumpl.Range("E11").Value = 7982.32
$wsCumpl.Range("F11").Value = 0.5056774832796631

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 22386.12
$wsCumpl.Range("E12").Value = 27920.88
$wsCumpl.Range("F12").Value = 0.4449901604150516

# Row 14 - TOTAL
$wsCumpl.Range("D14").Value = 60447.89999999999
$wsCumpl.Range("E14").Value = 37413.98766749098
$wsCumpl.Range("F14").Value = 0.6176858166213398

# Column E width tweak (24 -> 22 characters)
$wsCumpl.Columns.Item(5).ColumnWidth = 21.15
